# Update column C ("Förändrad") date value from 45224 to 45233 for rows 2-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45233
}
